# Applies the "Upload new version with timestamp" update to the DaySale report:
#  - adds a new item row "FLUB 20MG/ML SUSPENSION 30 ML" (inserted as new row 9,
#    pushing LASILACTONE..RICHI down by one row)
#  - adds a new item row "TIRATAM 100MG/ML ORAL SOLUTION 120 ML" as a new row
#    right after RICHI PANTHENOL (which is now row 13)
#  - updates RICHI PANTHENOL's stock/sales figures
#  - updates the grand total and the generated-at timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new data row at row 9 for "FLUB 20MG/ML SUSPENSION 30 ML".
#    This pushes the existing LASILACTONE / NEXICURE / OTRIVIN / RICHI rows
#    down from 9-12 to 10-13.
# ---------------------------------------------------------------------------
$ws.Rows.Item(9).Insert()

# Copy formatting from the row above (row 8) into the freshly inserted row 9,
# then fix up the row height + merges to match the other data rows.
$ws.Range("A8:Q8").Copy()
$ws.Range("A9:Q9").PasteSpecial(-4122)
$ws.Rows.Item(9).RowHeight = 25.5

$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "FLUB 20MG/ML SUSPENSION 30 ML"
$ws.Range("H9").Value = "3:0"
$ws.Range("L9").Value = 1
$ws.Range("N9").Value = "36.00"
$ws.Range("P9").Value = "72.0000"
$ws.Range("Q9").Value = "3:0"

# Renumber the following rows (A column counter 4..6 -> now at rows 10..12)
$ws.Range("A10").Value = 4
$ws.Range("A11").Value = 5
$ws.Range("A12").Value = 6

# ---------------------------------------------------------------------------
# 2) RICHI PANTHENOL ADVANCE GEL is now on row 13 - update its figures.
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = 7
$ws.Range("H13").Value = "0:0"
$ws.Range("P13").Value = "170.0000"
$ws.Range("Q13").Value = "3:0"

# ---------------------------------------------------------------------------
# 3) Insert a brand new row 14 for "TIRATAM 100MG/ML ORAL SOLUTION 120 ML"
#    right after RICHI PANTHENOL, before the totals/footer rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(14).Insert()

$ws.Range("A13:Q13").Copy()
$ws.Range("A14:Q14").PasteSpecial(-4122)
$ws.Rows.Item(14).RowHeight = 25.5

$ws.Range("A14:B14").Merge()
$ws.Range("C14:G14").Merge()
$ws.Range("H14:K14").Merge()
$ws.Range("L14:M14").Merge()
$ws.Range("N14:O14").Merge()

$ws.Range("A14").Value = 8
$ws.Range("C14").Value = "TIRATAM 100MG/ML ORAL SOLUTION 120 ML"
$ws.Range("H14").Value = "1:0"
$ws.Range("L14").Value = 1
$ws.Range("N14").Value = "120.00"
$ws.Range("P14").Value = "120.0000"
$ws.Range("Q14").Value = "1:0"

# ---------------------------------------------------------------------------
# 4) Update the grand total (now on row 15) and the generated-at timestamp
#    (now on row 16).
# ---------------------------------------------------------------------------
$ws.Range("P15").Value = 809
$ws.Range("A16").Value = "Tuesday, 19 August, 2025 10:04 AM"

Write-Output "edit applied"
